$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$donor = $ws.Range("B2").Style

$ws.Range("D2").Value = "'62.484.05"
$ws.Range("D2").Style = $donor
$ws.Range("E2").Value = "  +4.64%  "

$ws.Range("D3").Value = "'3.332.72"
$ws.Range("D3").Style = $donor
$ws.Range("E3").Value = "  +4.56%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'552.92"
$ws.Range("D5").Style = $donor
$ws.Range("E5").Value = "  +3.05%  "

$ws.Range("D6").Value = "'151.38"
$ws.Range("D6").Style = $donor
$ws.Range("E6").Value = "  +4.65%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").Style = $donor
$ws.Range("E8").Value = "  +2.20%  "

$ws.Range("D9").Value = "'7.51"
$ws.Range("D9").Style = $donor
$ws.Range("E9").Value = "  +3.08%  "

$ws.Range("D10").Value = "'0.117"
$ws.Range("D10").Style = $donor
$ws.Range("E10").Value = "  +3.85%  "

$ws.Range("E11").Value = "  +1.41%  "

$ws.Range("D12").Value = "'3.905.91"
$ws.Range("D12").Style = $donor
$ws.Range("E12").Value = "  +4.54%  "

$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "'0.0000180"
$ws.Range("D14").Style = $donor
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").Value = "'26.80"
$ws.Range("D15").Style = $donor
$ws.Range("E15").Value = "  +3.33%  "

$ws.Range("D16").Value = "'62.572.99"
$ws.Range("D16").Style = $donor
$ws.Range("E16").Value = "  +4.76%  "

$ws.Range("D17").Value = "'3.335.49"
$ws.Range("D17").Style = $donor
$ws.Range("E17").Value = "  +3.58%  "

$ws.Range("D18").Value = "'6.51"
$ws.Range("D18").Style = $donor
$ws.Range("E18").Value = "  +5.12%  "

$ws.Range("D19").Value = "'13.70"
$ws.Range("D19").Style = $donor
$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").Value = "'8.44"
$ws.Range("D20").Style = $donor
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").Value = "'385.58"
$ws.Range("D21").Style = $donor
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "'70.63"
$ws.Range("D24").Style = $donor
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("E25").Value = "  +2.95%  "

$ws.Range("D26").Value = "'8.90"
$ws.Range("D26").Style = $donor
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").Value = "'0.0₃0961"
$ws.Range("D27").Style = $donor
$ws.Range("E27").Value = "  +6.40%  "

$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = $donor
$ws.Range("E29").Value = "  +3.17%  "

$ws.Range("D30").Value = "'6.41"
$ws.Range("D30").Style = $donor
$ws.Range("E30").Value = "  +4.37%  "

$ws.Range("D31").Value = "'22.93"
$ws.Range("D31").Style = $donor
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("D32").Value = "'5.54"
$ws.Range("D32").Style = $donor
$ws.Range("E32").Value = "  +2.36%  "

$ws.Range("E33").Value = "  +8.16%  "

$ws.Range("D34").Value = "'6.72"
$ws.Range("D34").Style = $donor
$ws.Range("E34").Value = "  +4.03%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'161.40"
$ws.Range("D35").Style = $donor
$ws.Range("E35").Value = "  +2.99%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.48"
$ws.Range("D36").Style = $donor
$ws.Range("E36").Value = "  +9.89%  "

$ws.Range("D37").Value = "'1.86"
$ws.Range("D37").Style = $donor
$ws.Range("E37").Value = "  +10.51%  "

$ws.Range("E38").Value = "  +6.20%  "

$ws.Range("D39").Value = "'2.857.79"
$ws.Range("D39").Style = $donor
$ws.Range("E39").Value = "  +3.75%  "

$ws.Range("D40").Value = "'0.0740"
$ws.Range("D40").Style = $donor
$ws.Range("E40").Value = "  +3.94%  "

$ws.Range("D41").Value = "'0.0313"
$ws.Range("D41").Style = $donor
$ws.Range("E41").Value = "  +8.19%  "

$ws.Range("D42").Value = "'4.32"
$ws.Range("D42").Style = $donor
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").Value = "'0.751"
$ws.Range("D43").Style = $donor
$ws.Range("E43").Value = "  +3.52%  "

$ws.Range("D44").Value = "'40.49"
$ws.Range("D44").Style = $donor
$ws.Range("E44").Value = "  +2.61%  "

$ws.Range("E45").Value = "  +3.69%  "

$ws.Range("B46").Value = "RenzoRestakedETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D46").Value = "'3.378.63"
$ws.Range("D46").Style = $donor
$ws.Range("E46").Value = "  +4.49%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'21.92"
$ws.Range("D47").Style = $donor
$ws.Range("E47").Value = "  +6.89%  "

$ws.Range("E48").Value = "  +3.80%  "

$ws.Range("D49").Value = "'6.29"
$ws.Range("D49").Style = $donor
$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("D50").Value = "'0.803"
$ws.Range("D50").Style = $donor
$ws.Range("E50").Value = "  +3.29%  "

$ws.Range("D51").Value = "'282.50"
$ws.Range("D51").Style = $donor
$ws.Range("E51").Value = "  +8.63%  "
